$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 39 with the next forecast vector entry (ifo GDP component analysis preprocessing)
# Copy formatting from the last existing data row (A38) to the new date cell (A39)
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.2298740481777584
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = -0.05255865067609333
